# Regenerate the handback status report:
#   - first file: guid 46eaf1d9-a86d-4388-a250-c3f3bc71bcea -> 28d88ecf-cec1-483d-8822-62f010084577
#   - second file: guid fe9eac0d-dcbe-4fcd-b3f3-f790b30ae789 -> ffff41ab51ca-455e-4707-8e67-ba0ef629e9b3
#   - new xliff content hash 737ff6b4b7a31a98e040f8675f6154a9594cb646 / ad46878ae519ddeb0485bd1df7f3ee2c0aa183b1
#     collapse into a single regenerated hash 80c718f14bd74f96b3e433d7807a7ee00e2a0328
#   - refreshed handoff/handback timestamps

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$ws1 = $wb.Sheets.Item("Overview")

$ws1.Range("A2").Value = "28d88ecf-cec1-483d-8822-62f010084577.md"
$ws1.Range("B2").Value = "e2e\28d88ecf-cec1-483d-8822-62f010084577.md"
$ws1.Range("G2").Value = "2016-09-04 03:07:10"

$ws1.Range("A3").Value = "ffff41ab51ca-455e-4707-8e67-ba0ef629e9b3.md"
$ws1.Range("B3").Value = "e2e\ffff41ab51ca-455e-4707-8e67-ba0ef629e9b3.md"
$ws1.Range("G3").Value = "2016-09-04 03:07:10"

# Rebuild the two B-column hyperlinks so the display text tracks the new file
# names while the link targets stay the same as before.
$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/679a39655bc1552ebfa7b31f609ef3bb074eaabe/e2e/46eaf1d9-a86d-4388-a250-c3f3bc71bcea.md", [Type]::Missing, [Type]::Missing, "e2e\28d88ecf-cec1-483d-8822-62f010084577.md")
$ws1.Hyperlinks.Add($ws1.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/679a39655bc1552ebfa7b31f609ef3bb074eaabe/e2e/fe9eac0d-dcbe-4fcd-b3f3-f790b30ae789.md", [Type]::Missing, [Type]::Missing, "e2e\ffff41ab51ca-455e-4707-8e67-ba0ef629e9b3.md")

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$ws2 = $wb.Sheets.Item("zh-cn")

$ws2.Range("A2").Value = "28d88ecf-cec1-483d-8822-62f010084577.md"
$ws2.Range("G2").Value = "28d88ecf-cec1-483d-8822-62f010084577.80c718f14bd74f96b3e433d7807a7ee00e2a0328.zh-cn.xlf"
$ws2.Range("H2").Value = "2016-09-04 03:07:01"
$ws2.Range("I2").Value = "28d88ecf-cec1-483d-8822-62f010084577.md"
$ws2.Range("J2").Value = "28d88ecf-cec1-483d-8822-62f010084577.80c718f14bd74f96b3e433d7807a7ee00e2a0328.zh-cn.xlf"
$ws2.Range("K2").Value = "2016-09-04 03:07:51"

$ws2.Range("A3").Value = "ffff41ab51ca-455e-4707-8e67-ba0ef629e9b3.md"
$ws2.Range("G3").Value = "28d88ecf-cec1-483d-8822-62f010084577.80c718f14bd74f96b3e433d7807a7ee00e2a0328.zh-cn.xlf"
$ws2.Range("H3").Value = "2016-09-04 03:07:01"
$ws2.Range("I3").Value = "ffff41ab51ca-455e-4707-8e67-ba0ef629e9b3.md"
$ws2.Range("J3").Value = "28d88ecf-cec1-483d-8822-62f010084577.80c718f14bd74f96b3e433d7807a7ee00e2a0328.zh-cn.xlf"
$ws2.Range("K3").Value = "2016-09-04 03:07:51"

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/679a39655bc1552ebfa7b31f609ef3bb074eaabe/e2e/46eaf1d9-a86d-4388-a250-c3f3bc71bcea.md", [Type]::Missing, [Type]::Missing, "28d88ecf-cec1-483d-8822-62f010084577.md")
$ws2.Hyperlinks.Add($ws2.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/bcce357308de4dd1babfcd3a5c3fb38cee325f80/e2e/46eaf1d9-a86d-4388-a250-c3f3bc71bcea.md", [Type]::Missing, [Type]::Missing, "28d88ecf-cec1-483d-8822-62f010084577.md")
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/679a39655bc1552ebfa7b31f609ef3bb074eaabe/e2e/fe9eac0d-dcbe-4fcd-b3f3-f790b30ae789.md", [Type]::Missing, [Type]::Missing, "ffff41ab51ca-455e-4707-8e67-ba0ef629e9b3.md")
$ws2.Hyperlinks.Add($ws2.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/bcce357308de4dd1babfcd3a5c3fb38cee325f80/e2e/fe9eac0d-dcbe-4fcd-b3f3-f790b30ae789.md", [Type]::Missing, [Type]::Missing, "ffff41ab51ca-455e-4707-8e67-ba0ef629e9b3.md")

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$ws3 = $wb.Sheets.Item("de-de")

$ws3.Range("A2").Value = "28d88ecf-cec1-483d-8822-62f010084577.md"
$ws3.Range("G2").Value = "28d88ecf-cec1-483d-8822-62f010084577.80c718f14bd74f96b3e433d7807a7ee00e2a0328.de-de.xlf"
$ws3.Range("H2").Value = "2016-09-04 03:07:10"
$ws3.Range("I2").Value = "28d88ecf-cec1-483d-8822-62f010084577.md"
$ws3.Range("J2").Value = "28d88ecf-cec1-483d-8822-62f010084577.80c718f14bd74f96b3e433d7807a7ee00e2a0328.de-de.xlf"
$ws3.Range("K2").Value = "2016-09-04 03:07:59"

$ws3.Range("A3").Value = "ffff41ab51ca-455e-4707-8e67-ba0ef629e9b3.md"
$ws3.Range("G3").Value = "28d88ecf-cec1-483d-8822-62f010084577.80c718f14bd74f96b3e433d7807a7ee00e2a0328.de-de.xlf"
$ws3.Range("H3").Value = "2016-09-04 03:07:10"
$ws3.Range("I3").Value = "ffff41ab51ca-455e-4707-8e67-ba0ef629e9b3.md"
$ws3.Range("J3").Value = "28d88ecf-cec1-483d-8822-62f010084577.80c718f14bd74f96b3e433d7807a7ee00e2a0328.de-de.xlf"
$ws3.Range("K3").Value = "2016-09-04 03:07:59"

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/679a39655bc1552ebfa7b31f609ef3bb074eaabe/e2e/46eaf1d9-a86d-4388-a250-c3f3bc71bcea.md", [Type]::Missing, [Type]::Missing, "28d88ecf-cec1-483d-8822-62f010084577.md")
$ws3.Hyperlinks.Add($ws3.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/4534189d8b520a25e172aa99a14ce3ccfd60abdd/e2e/46eaf1d9-a86d-4388-a250-c3f3bc71bcea.md", [Type]::Missing, [Type]::Missing, "28d88ecf-cec1-483d-8822-62f010084577.md")
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/679a39655bc1552ebfa7b31f609ef3bb074eaabe/e2e/fe9eac0d-dcbe-4fcd-b3f3-f790b30ae789.md", [Type]::Missing, [Type]::Missing, "ffff41ab51ca-455e-4707-8e67-ba0ef629e9b3.md")
$ws3.Hyperlinks.Add($ws3.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/4534189d8b520a25e172aa99a14ce3ccfd60abdd/e2e/fe9eac0d-dcbe-4fcd-b3f3-f790b30ae789.md", [Type]::Missing, [Type]::Missing, "ffff41ab51ca-455e-4707-8e67-ba0ef629e9b3.md")
